# Apply "notes of references, stuggle with bib" changes to the
# "Interference sources" sheet, and move the active tab back to it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Interference sources")

# --- New "Mech-Interference" section ------------------------------------
$ws.Range("A29").Value = "Mech-Interference"

$ws.Range("A30").Value = "hyper-CO interference"
$ws.Range("A30").Font.Bold = $true

# --- Row 26: mark existing "Otto Payseur" entry bold and add its
#     associated notes / title in C26 / D26 -----------------------------
$ws.Range("B26").Font.Bold = $true
$ws.Range("D26").Value = "Crossover Interference: Shedding Light on the Evolution of Recombination"
$ws.Range("C26").Value = "theo, simulations of interference variation"

$ws.Range("A32").Value = "natural variation in COI"
$ws.Range("A32").Font.Bold = $true

$ws.Range("A33").Value = 2013
$ws.Range("A34").Value = 2016
$ws.Range("B34").Value = "Wang Z, Shen B, Jiang J, Li J, Ma L. 2016. Effect of sex, age and genetics on crossover interference in cattle. Sci. Rep. 6:37698"
$ws.Range("B33").Value = "Bauer E, Falque M,Walter H, Bauland C, Camisan C, et al. 2013. Intraspecific variation of recombination rate in maize. Genome Biol. 14:R103–1"

# --- View state: move selection / frozen-pane scroll, and switch the
#     active tab back to "Interference sources" -------------------------
$ws2 = $wb.Worksheets.Item("DSB.DMC1")
[void]$ws2.Activate()
$ws2.Range("B33").Select() | Out-Null

[void]$ws.Activate()
$ws.Range("A11").Select() | Out-Null
$ws.Range("C28").Select() | Out-Null
